# Inventario: carga masiva de ventas de enero, ajuste de cantidades
# existentes y nuevas filas de "Ingresos" (Barra / Cafeteria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Actualizacion masiva de la columna "cantidad" (C) ---------------
$qtyUpdates = @{
    "C2"  = 17
    "C3"  = 47
    "C4"  = 12
    "C5"  = 15
    "C6"  = 25
    "C7"  = 20
    "C8"  = 11
    "C9"  = 0
    "C10" = 31
    "C11" = 27
    "C13" = 24
    "C15" = 27
    "C16" = 19
    "C17" = 25
    "C18" = 12
    "C20" = 18
    "C23" = 144
    "C24" = 117
    "C25" = 2
    "C26" = 97
    "C29" = 39
    "C30" = 5
    "C31" = 42
    "C32" = 39
    "C33" = 45
    "C34" = 38
    "C37" = 28
    "C38" = 21
    "C39" = 7
}

foreach ($addr in $qtyUpdates.Keys) {
    $ws.Range($addr).Value = $qtyUpdates[$addr]
}

# --- Nuevas filas: ingresos de Barra y Cafeteria ----------------------
$newRows = @(
    @{ row = 42; nombre = "Ventas Barra";      categoria = "Ingresos Barra";      cantidad = 1000; precio = 0; fecha = "2/22/2026"; costo = 0 },
    @{ row = 43; nombre = "Ventas Cafeteria";   categoria = "Ingresos Cafeteria";  cantidad = 1000; precio = 0; fecha = "2/22/2026"; costo = 0 }
)

foreach ($r in $newRows) {
    $n = $r.row
    $ws.Range("A$n").Value = $r.nombre
    $ws.Range("B$n").Value = $r.categoria
    $ws.Range("C$n").Value = $r.cantidad
    $ws.Range("D$n").Value = $r.precio
    # Prefijo de comilla para forzar texto literal (evita que "2/22/2026"
    # se autoconvierta a un valor de fecha), igual que en las filas previas.
    $ws.Range("E$n").Value = "'" + $r.fecha
    $ws.Range("F$n").Value = $r.costo
}
